$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get reshuffled between rows (D, H, I, J, K, L, M, N, P, Q).
# A, B, C, E, F, G, O, R stay identical for every row (market/region/category
# metadata), so only these columns need to move.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "P", "Q")

# Mapping of destination row -> source row (taken from the source row's
# original values before any writes happen), derived from the diff: this is
# a weekly re-shuffle of the price records among the existing dated rows.
$rowMap = @{
    2  = 15
    3  = 11
    4  = 4
    5  = 7
    6  = 8
    7  = 9
    8  = 12
    9  = 6
    10 = 5
    11 = 16
    12 = 17
    13 = 2
    14 = 10
    15 = 18
    16 = 14
    17 = 3
    18 = 13
}

# Snapshot every source row's values first, so overwriting a row doesn't
# clobber data that still needs to be read for another destination row.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$srcRow").Value2
        }
        $snapshot[$srcRow] = $rowData
    }
}

# Now write the snapshot values into their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$col]
    }
}
